$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.80"
$ws.Range("E2").Value = "'1.07%"

$ws.Range("D3").Value = "'43.98"
$ws.Range("E3").Value = "'-1.10%"

$ws.Range("D4").Value = "'5.500"
$ws.Range("E4").Value = "'0.12%"

$ws.Range("E5").Value = "'-0.64%"

$ws.Range("D6").Value = "'1.993"
$ws.Range("E6").Value = "'4.59%"

$ws.Range("B7").Value = 'BTSEToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D7").Value = "'2.574"
$ws.Range("E7").Value = "'-4.82%"

$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = "'0.9493"
$ws.Range("E8").Value = "'0.90%"

$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = "'0.1125"
$ws.Range("E9").Value = "'-3.93%"

$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = "'0.1888"
$ws.Range("E10").Value = "'1.10%"

$ws.Range("B11").Value = 'MCDex'
$ws.Range("C11").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D11").Value = "'10.63"
$ws.Range("E11").Value = "'25.82%"

$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = "'0.1008"
$ws.Range("E12").Value = "'0.24%"

$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = "'0.04792"
$ws.Range("E13").Value = "'12.53%"

$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = "'0.1062"
$ws.Range("E14").Value = "'-0.20%"

$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = "'0.001274"
$ws.Range("E15").Value = "'-0.29%"

$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").Value = "'0.04076"
$ws.Range("E16").Value = "'-3.80%"

$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = "'0.005990"
$ws.Range("E17").Value = "'2.03%"

$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = "'3.364"
$ws.Range("E18").Value = "'-6.14%"

$ws.Range("B19").Value = 'GateToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D19").Value = "'4.375"
$ws.Range("E19").Value = "'2.32%"

$ws.Range("E20").Value = "'-0.68%"

$ws.Range("E21").Value = "'3.04%"

$ws.Range("E23").Value = "'2.45%"

$ws.Range("D24").Value = "'0.004334"
$ws.Range("E24").Value = "'-4.87%"

$ws.Range("D25").Value = "'0.0001200"
$ws.Range("E25").Value = "'1.68%"

$ws.Range("E26").Value = "'-6.15%"

$ws.Range("D38").Value = "'0.02583"
$ws.Range("E38").Value = "'-2.06%"

$ws.Range("D39").Value = "'0.05658"
$ws.Range("E39").Value = "'3.28%"

$ws.Range("D40").Value = "'0.007547"
$ws.Range("E40").Value = "'-1.61%"

$ws.Range("E41").Value = "'0.02%"

$ws.Range("D42").Value = "'0.007409"
$ws.Range("E42").Value = "'3.52%"

$ws.Range("D43").Value = "'0.002015"
$ws.Range("E43").Value = "'-0.04%"

$ws.Range("D44").Value = "'0.008635"
$ws.Range("E44").Value = "'-6.12%"

$ws.Range("D45").Value = "'0.00007109"
$ws.Range("E45").Value = "'0.22%"

$ws.Range("E46").Value = "'-0.01%"

$ws.Range("E47").Value = "'55.48%"

$ws.Range("D48").Value = "'0.003789"
$ws.Range("E48").Value = "'5.66%"

$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'-0.01%"

$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'-0.01%"
